# shorter elimination group labels
#
# - adds a new "match.form.label-inverse" row to the Fonts table (Table1 on
#   the "Fonts" sheet), supplying an italic Consola variant for the inverse
#   match-form label
# - repoints "match.form.label" at the regular Consola font (en/ja/fa)
#   instead of the old TradeGothic/YuGoth fonts
# - re-sorts the Fonts table alphabetically by its "key" column
# - updates the view (active sheet/selection/zoom/column width) to match
#   where the author ended up

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fonts")
$tbl = $ws.ListObjects.Item(1)

# --- update the existing "match.form.label" row to use Consola fonts ---
for ($i = 1; $i -le $tbl.ListRows.Count; $i++) {
    $keyCell = $ws.Cells.Item($i + 1, 1)
    if ($keyCell.Value2 -eq "match.form.label") {
        $ws.Cells.Item($i + 1, 2).Value = "consola.ttf"
        $ws.Cells.Item($i + 1, 9).Value = "consola.ttf"
        $ws.Cells.Item($i + 1, 10).Value = "consola.ttf"
    }
}

# --- add the new "match.form.label-inverse" row ---
$newRow = $tbl.ListRows.Add()
$newIdx = $newRow.Index
$ws.Cells.Item($newIdx + 1, 1).Value = "match.form.label-inverse"
$ws.Cells.Item($newIdx + 1, 2).Value = "consolai.ttf"
$ws.Cells.Item($newIdx + 1, 3).Value = "calibri.ttf"
$ws.Cells.Item($newIdx + 1, 9).Value = "consolai.ttf"
$ws.Cells.Item($newIdx + 1, 10).Value = "consolai.ttf"

# --- re-sort the table alphabetically by key ---
$tbl.Sort.SortFields.Clear()
$tbl.Sort.SortFields.Add($tbl.ListColumns.Item(1).Range, 0, 1)
$tbl.Sort.Header = 1
$tbl.Sort.Apply()

# --- widen column A to fit the longer keys ---
$ws.Columns.Item(1).ColumnWidth = 19.669256718124647

# --- switch to the Fonts sheet, set zoom + selection to match the author ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 150
$ws.Range("I16:I17").Select()
